{"js": "// Update the worksheet date heading and the 25 \"two-digit \u00f7 one-digit\"\n// division problems laid out in the document's 5x5 problem grid.\n//\n// Each replacement is looked up by its exact, unique original text via\n// Body.search() and the matching run is replaced in place, which keeps\n// the existing run/paragraph formatting (font, size, alignment) intact.\n\nconst body = context.document.body;\n\nasync function replaceOnce(oldText, newText) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(newText, Word.InsertLocation.replace);\n    await context.sync();\n  }\n}\n\n// 1) Date / weekday heading.\nawait replaceOnce(\"2024-07-19 Friday\", \"2024-07-20 Saturday\");\n\n// 2) Division problems. Listed in document (reading) order. The last two\n// entries are intentionally swapped relative to reading order (\"81\u00f79=\"\n// is handled before \"11\u00f74=\") because \"11\u00f74=\" is replaced with \"81\u00f79=\",\n// which would otherwise collide with the still-unprocessed original\n// \"81\u00f79=\" cell.\nconst problemReplacements = [\n  [\"26\u00f73=\", \"18\u00f77=\"],\n  [\"77\u00f73=\", \"21\u00f78=\"],\n  [\"35\u00f79=\", \"15\u00f73=\"],\n  [\"33\u00f79=\", \"79\u00f76=\"],\n  [\"61\u00f72=\", \"36\u00f79=\"],\n  [\"42\u00f76=\", \"36\u00f79=\"],\n  [\"13\u00f72=\", \"73\u00f75=\"],\n  [\"26\u00f79=\", \"39\u00f79=\"],\n  [\"65\u00f72=\", \"80\u00f75=\"],\n  [\"73\u00f79=\", \"36\u00f77=\"],\n  [\"49\u00f77=\", \"39\u00f76=\"],\n  [\"53\u00f74=\", \"71\u00f73=\"],\n  [\"47\u00f75=\", \"86\u00f74=\"],\n  [\"14\u00f75=\", \"55\u00f73=\"],\n  [\"49\u00f72=\", \"85\u00f78=\"],\n  [\"99\u00f72=\", \"40\u00f73=\"],\n  [\"29\u00f78=\", \"78\u00f79=\"],\n  [\"42\u00f75=\", \"87\u00f73=\"],\n  [\"98\u00f77=\", \"72\u00f74=\"],\n  [\"28\u00f79=\", \"60\u00f75=\"],\n  [\"59\u00f77=\", \"61\u00f75=\"],\n  [\"87\u00f78=\", \"68\u00f74=\"],\n  [\"52\u00f78=\", \"29\u00f72=\"],\n  [\"81\u00f79=\", \"67\u00f73=\"],\n  [\"11\u00f74=\", \"81\u00f79=\"],\n];\n\nfor (const [oldText, newText] of problemReplacements) {\n  await replaceOnce(oldText, newText);\n}\n", "ps1": "# Update the worksheet date heading and the 25 \"two-digit / one-digit\"\n# division problems laid out in the document's 5x5 problem grid.\n#\n# Each replacement targets its exact, unique original text with\n# Find/Replace restricted to a single occurrence (wdReplaceOne), which\n# leaves the surrounding run/paragraph formatting (font, size,\n# alignment) untouched.\n\n$d = $word.ActiveDocument\n\nfunction Replace-Once {\n    param(\n        [string]$Find,\n        [string]$Replace\n    )\n\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $range.Find.Text = $Find\n    $range.Find.Replacement.Text = $Replace\n    # MatchCase:=True, Forward:=True, Wrap:=wdFindContinue(1), Replace:=wdReplaceOne(1)\n    $range.Find.Execute($Find, $true, $false, $false, $false, $false, $true, 1, $false, $Replace, 1) | Out-Null\n}\n\n# 1) Date / weekday heading.\nReplace-Once \"2024-07-19 Friday\" \"2024-07-20 Saturday\"\n\n# 2) Division problems, in document (reading) order. The last two calls\n# are intentionally swapped relative to reading order (\"81\u00f79=\" is\n# handled before \"11\u00f74=\") because \"11\u00f74=\" is replaced with \"81\u00f79=\",\n# which would otherwise collide with the still-unprocessed original\n# \"81\u00f79=\" cell.\nReplace-Once \"26\u00f73=\" \"18\u00f77=\"\nReplace-Once \"77\u00f73=\" \"21\u00f78=\"\nReplace-Once \"35\u00f79=\" \"15\u00f73=\"\nReplace-Once \"33\u00f79=\" \"79\u00f76=\"\nReplace-Once \"61\u00f72=\" \"36\u00f79=\"\nReplace-Once \"42\u00f76=\" \"36\u00f79=\"\nReplace-Once \"13\u00f72=\" \"73\u00f75=\"\nReplace-Once \"26\u00f79=\" \"39\u00f79=\"\nReplace-Once \"65\u00f72=\" \"80\u00f75=\"\nReplace-Once \"73\u00f79=\" \"36\u00f77=\"\nReplace-Once \"49\u00f77=\" \"39\u00f76=\"\nReplace-Once \"53\u00f74=\" \"71\u00f73=\"\nReplace-Once \"47\u00f75=\" \"86\u00f74=\"\nReplace-Once \"14\u00f75=\" \"55\u00f73=\"\nReplace-Once \"49\u00f72=\" \"85\u00f78=\"\nReplace-Once \"99\u00f72=\" \"40\u00f73=\"\nReplace-Once \"29\u00f78=\" \"78\u00f79=\"\nReplace-Once \"42\u00f75=\" \"87\u00f73=\"\nReplace-Once \"98\u00f77=\" \"72\u00f74=\"\nReplace-Once \"28\u00f79=\" \"60\u00f75=\"\nReplace-Once \"59\u00f77=\" \"61\u00f75=\"\nReplace-Once \"87\u00f78=\" \"68\u00f74=\"\nReplace-Once \"52\u00f78=\" \"29\u00f72=\"\nReplace-Once \"81\u00f79=\" \"67\u00f73=\"\nReplace-Once \"11\u00f74=\" \"81\u00f79=\"\n"}
